{"js": "// Replace the 25 \"three-digit x one-digit\" practice answers in the\n// table with their updated values. Every source string is unique within\n// the document, so a targeted body.search() + insertText(replace) per\n// pair is unambiguous and leaves every other run property untouched.\nconst replacements = [\n  [\"614\u00d76=3684\", \"259\u00d75=1295\"],\n  [\"713\u00d73=2139\", \"564\u00d72=1128\"],\n  [\"133\u00d78=1064\", \"516\u00d73=1548\"],\n  [\"342\u00d72=684\", \"862\u00d73=2586\"],\n  [\"781\u00d79=7029\", \"426\u00d77=2982\"],\n  [\"124\u00d77=868\", \"646\u00d76=3876\"],\n  [\"723\u00d73=2169\", \"449\u00d73=1347\"],\n  [\"309\u00d76=1854\", \"610\u00d78=4880\"],\n  [\"249\u00d73=747\", \"780\u00d72=1560\"],\n  [\"898\u00d73=2694\", \"722\u00d79=6498\"],\n  [\"749\u00d75=3745\", \"388\u00d75=1940\"],\n  [\"271\u00d78=2168\", \"526\u00d77=3682\"],\n  [\"664\u00d76=3984\", \"981\u00d78=7848\"],\n  [\"884\u00d73=2652\", \"693\u00d79=6237\"],\n  [\"404\u00d72=808\", \"936\u00d75=4680\"],\n  [\"540\u00d77=3780\", \"701\u00d76=4206\"],\n  [\"801\u00d73=2403\", \"388\u00d73=1164\"],\n  [\"796\u00d76=4776\", \"122\u00d79=1098\"],\n  [\"935\u00d72=1870\", \"780\u00d76=4680\"],\n  [\"498\u00d72=996\", \"804\u00d77=5628\"],\n  [\"301\u00d74=1204\", \"917\u00d77=6419\"],\n  [\"164\u00d79=1476\", \"419\u00d72=838\"],\n  [\"838\u00d78=6704\", \"141\u00d75=705\"],\n  [\"671\u00d79=6039\", \"963\u00d73=2889\"],\n  [\"350\u00d77=2450\", \"843\u00d75=4215\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"three-digit x one-digit\" practice answers in the table\n# with their updated values. Every source string is unique within the\n# document, so Find/Replace against $d.Content (whole-document range) for\n# each pair is unambiguous and leaves every other run property untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{ old = \"614\u00d76=3684\"; new = \"259\u00d75=1295\" },\n  @{ old = \"713\u00d73=2139\"; new = \"564\u00d72=1128\" },\n  @{ old = \"133\u00d78=1064\"; new = \"516\u00d73=1548\" },\n  @{ old = \"342\u00d72=684\"; new = \"862\u00d73=2586\" },\n  @{ old = \"781\u00d79=7029\"; new = \"426\u00d77=2982\" },\n  @{ old = \"124\u00d77=868\"; new = \"646\u00d76=3876\" },\n  @{ old = \"723\u00d73=2169\"; new = \"449\u00d73=1347\" },\n  @{ old = \"309\u00d76=1854\"; new = \"610\u00d78=4880\" },\n  @{ old = \"249\u00d73=747\"; new = \"780\u00d72=1560\" },\n  @{ old = \"898\u00d73=2694\"; new = \"722\u00d79=6498\" },\n  @{ old = \"749\u00d75=3745\"; new = \"388\u00d75=1940\" },\n  @{ old = \"271\u00d78=2168\"; new = \"526\u00d77=3682\" },\n  @{ old = \"664\u00d76=3984\"; new = \"981\u00d78=7848\" },\n  @{ old = \"884\u00d73=2652\"; new = \"693\u00d79=6237\" },\n  @{ old = \"404\u00d72=808\"; new = \"936\u00d75=4680\" },\n  @{ old = \"540\u00d77=3780\"; new = \"701\u00d76=4206\" },\n  @{ old = \"801\u00d73=2403\"; new = \"388\u00d73=1164\" },\n  @{ old = \"796\u00d76=4776\"; new = \"122\u00d79=1098\" },\n  @{ old = \"935\u00d72=1870\"; new = \"780\u00d76=4680\" },\n  @{ old = \"498\u00d72=996\"; new = \"804\u00d77=5628\" },\n  @{ old = \"301\u00d74=1204\"; new = \"917\u00d77=6419\" },\n  @{ old = \"164\u00d79=1476\"; new = \"419\u00d72=838\" },\n  @{ old = \"838\u00d78=6704\"; new = \"141\u00d75=705\" },\n  @{ old = \"671\u00d79=6039\"; new = \"963\u00d73=2889\" },\n  @{ old = \"350\u00d77=2450\"; new = \"843\u00d75=4215\" }\n)\n\nforeach ($pair in $pairs) {\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n  if (-not $found) {\n    throw \"Text not found: $($pair.old)\"\n  }\n}\n"}
